# Auto-generated market data refresh for Pandaemonium_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) per leve row
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2728.0715
$ws.Range("I40").Value = 2917.5454
$ws.Range("J40").Value = 2033.3334
$ws.Range("K40").Value = 2917.5454
$ws.Range("L40").Value = 2033.3334
$ws.Range("M40").Value = -2742.5454
$ws.Range("N40").Value = -2383.3334

$ws.Range("H69").Value = 7004.3335
$ws.Range("I69").Value = 8006.5
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 24019.5
$ws.Range("L69").Value = 15000
$ws.Range("M69").Value = -23145.5
$ws.Range("N69").Value = -16748

$ws.Range("H72").Value = 7004.3335
$ws.Range("I72").Value = 8006.5
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 72058.5
$ws.Range("L72").Value = 45000
$ws.Range("M72").Value = -67690.5
$ws.Range("N72").Value = -53736

$ws.Range("H108").Value = 42776.5
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 42776.5
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 42776.5
$ws.Range("N108").Value = -50456.5

$ws.Range("H109").Value = 59596.668
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 59596.668
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 59596.668
$ws.Range("N109").Value = -62370.668

$ws.Range("H110").Value = 36600.43
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 36600.43
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 36600.43
$ws.Range("N110").Value = -44780.43

$ws.Range("H125").Value = 7910.8667
$ws.Range("I125").Value = 1714.2858
$ws.Range("J125").Value = 13332.875
$ws.Range("K125").Value = 15428.5722
$ws.Range("L125").Value = 119995.875
$ws.Range("M125").Value = -12968.5722
$ws.Range("N125").Value = -124915.875

$ws.Range("H132").Value = 1735.4667
$ws.Range("I132").Value = 1986.4828
$ws.Range("J132").Value = 1280.5
$ws.Range("K132").Value = 5959.4484
$ws.Range("L132").Value = 3841.5
$ws.Range("M132").Value = -3429.4484
$ws.Range("N132").Value = -8901.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 2000
$ws.Range("I12").Value = 2000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 2000
$ws.Range("L12").ClearContents()
$ws.Range("M12").Value = -1827
$ws.Range("N12").Value = 0

$ws.Range("H32").Value = 5799.25
$ws.Range("I32").Value = 5065.894
$ws.Range("J32").Value = 30000
$ws.Range("K32").Value = 5065.894
$ws.Range("L32").Value = 30000
$ws.Range("M32").Value = -4778.894
$ws.Range("N32").Value = -30574

$ws.Range("H61").Value = 6234.5713
$ws.Range("I61").Value = 5675.1333
$ws.Range("J61").Value = 7633.1665
$ws.Range("K61").Value = 5675.1333
$ws.Range("L61").Value = 7633.1665
$ws.Range("M61").Value = -5463.1333
$ws.Range("N61").Value = -8057.1665

$ws.Range("H88").Value = 5216.4
$ws.Range("I88").Value = 7876.5
$ws.Range("J88").Value = 2176.2856
$ws.Range("K88").Value = 7876.5
$ws.Range("L88").Value = 2176.2856
$ws.Range("M88").Value = -7470.5
$ws.Range("N88").Value = -2988.2856

$ws.Range("H91").Value = 5216.4
$ws.Range("I91").Value = 7876.5
$ws.Range("J91").Value = 2176.2856
$ws.Range("K91").Value = 7876.5
$ws.Range("L91").Value = 2176.2856
$ws.Range("M91").Value = -6472.5
$ws.Range("N91").Value = -4984.2856

$ws.Range("H122").Value = 3270.4666
$ws.Range("I122").Value = 3727.3333
$ws.Range("J122").Value = 2585.1667
$ws.Range("K122").Value = 11181.9999
$ws.Range("L122").Value = 7755.500100000001
$ws.Range("M122").Value = -8731.999899999999
$ws.Range("N122").Value = -12655.5001

$ws.Range("H136").Value = 6234.5713
$ws.Range("I136").Value = 5675.1333
$ws.Range("J136").Value = 7633.1665
$ws.Range("K136").Value = 17025.3999
$ws.Range("L136").Value = 22899.4995
$ws.Range("M136").Value = -14475.3999
$ws.Range("N136").Value = -27999.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 383.33334
$ws.Range("I11").Value = 383.33334
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 383.33334
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -243.33334

$ws.Range("H134").Value = 3911.2285
$ws.Range("I134").Value = 3285.9656
$ws.Range("J134").Value = 6933.3335
$ws.Range("K134").Value = 9857.8968
$ws.Range("L134").Value = 20800.0005
$ws.Range("M134").Value = -7322.8968
$ws.Range("N134").Value = -25870.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1983.0483
$ws.Range("I31").Value = 1618.3334
$ws.Range("J31").Value = 2601.4783
$ws.Range("K31").Value = 1618.3334
$ws.Range("L31").Value = 2601.4783
$ws.Range("M31").Value = -1323.3334
$ws.Range("N31").Value = -3191.4783

$ws.Range("H34").Value = 1983.0483
$ws.Range("I34").Value = 1618.3334
$ws.Range("J34").Value = 2601.4783
$ws.Range("K34").Value = 1618.3334
$ws.Range("L34").Value = 2601.4783
$ws.Range("M34").Value = -1416.3334
$ws.Range("N34").Value = -3005.4783

$ws.Range("H105").Value = 1067.3572
$ws.Range("I105").Value = 914.8
$ws.Range("J105").Value = 1448.75
$ws.Range("K105").Value = 914.8
$ws.Range("L105").Value = 1448.75
$ws.Range("M105").Value = 832.2
$ws.Range("N105").Value = -4942.75

$ws.Range("H132").Value = 2571.5
$ws.Range("I132").Value = 1949.4445
$ws.Range("J132").Value = 3691.2
$ws.Range("K132").Value = 5848.333500000001
$ws.Range("L132").Value = 11073.6
$ws.Range("M132").Value = -3318.333500000001
$ws.Range("N132").Value = -16133.6

$ws.Range("H134").Value = 2889.2856
$ws.Range("I134").Value = 1759.4814
$ws.Range("J134").Value = 3941.1724
$ws.Range("K134").Value = 5278.4442
$ws.Range("L134").Value = 11823.5172
$ws.Range("M134").Value = -2743.4442
$ws.Range("N134").Value = -16893.5172

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2389.1072
$ws.Range("I34").Value = 220.5
$ws.Range("J34").Value = 3593.889
$ws.Range("K34").Value = 661.5
$ws.Range("L34").Value = 10781.667
$ws.Range("M34").Value = -577.5
$ws.Range("N34").Value = -10949.667

$ws.Range("H46").Value = 2477.7778
$ws.Range("I46").Value = 200
$ws.Range("J46").Value = 3128.5715
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 9385.7145
$ws.Range("M46").Value = -509
$ws.Range("N46").Value = -9567.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10329.286
$ws.Range("I80").Value = 26152.5
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 26152.5
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -25154.5
$ws.Range("N80").Value = -5996

$ws.Range("H83").Value = 10329.286
$ws.Range("I83").Value = 26152.5
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 130762.5
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -125770.5
$ws.Range("N83").Value = -29984

$ws.Range("H132").Value = 2471.5588
$ws.Range("I132").Value = 2275.762
$ws.Range("J132").Value = 2787.8462
$ws.Range("K132").Value = 6827.286
$ws.Range("L132").Value = 8363.5386
$ws.Range("M132").Value = -4297.286
$ws.Range("N132").Value = -13423.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 37975
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 37975
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 37975
$ws.Range("N110").Value = -46155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3291.6667
$ws.Range("I122").Value = 1780.9375
$ws.Range("J122").Value = 6313.125
$ws.Range("K122").Value = 5342.8125
$ws.Range("L122").Value = 18939.375
$ws.Range("M122").Value = -2892.8125
$ws.Range("N122").Value = -23839.375

$ws.Range("H126").Value = 1358.4
$ws.Range("I126").Value = 845.8
$ws.Range("J126").Value = 1871
$ws.Range("K126").Value = 2537.4
$ws.Range("L126").Value = 5613
$ws.Range("M126").Value = -67.39999999999964
$ws.Range("N126").Value = -10553

$ws.Range("H132").Value = 1285.2894
$ws.Range("I132").Value = 1183.3667
$ws.Range("J132").Value = 1667.5
$ws.Range("K132").Value = 3550.1001
$ws.Range("L132").Value = 5002.5
$ws.Range("M132").Value = -1020.1001
$ws.Range("N132").Value = -10062.5

$ws.Range("H136").Value = 3725.5
$ws.Range("I136").Value = 1799.8889
$ws.Range("J136").Value = 5764.3823
$ws.Range("K136").Value = 5399.6667
$ws.Range("L136").Value = 17293.1469
$ws.Range("M136").Value = -2849.6667
$ws.Range("N136").Value = -22393.1469

